$d = $word.ActiveDocument

# Update the date heading.
$d.Content.Find.Execute("2023-12-16 Saturday", $false, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-17 Sunday", 2)

# Update the multiplication problems in the table, cell by cell
# (some old values repeat, e.g. "53x84=", so address cells directly
# rather than doing a blanket find/replace).
$t = $d.Tables.Item(1)

$values = @(
    @(1,1,"66×72="),
    @(1,2,"22×90="),
    @(1,3,"89×36="),
    @(1,4,"72×96="),
    @(1,5,"79×69="),

    @(5,1,"17×74="),
    @(5,2,"19×37="),
    @(5,3,"13×35="),
    @(5,4,"41×23="),
    @(5,5,"62×77="),

    @(10,1,"45×23="),
    @(10,2,"15×44="),
    @(10,3,"13×93="),
    @(10,4,"80×18="),
    @(10,5,"60×59="),

    @(15,1,"64×34="),
    @(15,2,"73×81="),
    @(15,3,"24×20="),
    @(15,4,"69×98="),
    @(15,5,"36×64="),

    @(20,1,"51×34="),
    @(20,2,"24×40="),
    @(20,3,"64×79="),
    @(20,4,"53×56="),
    @(20,5,"18×66=")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $col = $entry[1]
    $newText = $entry[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
